$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-8 (project IDs 5-7) change Status from RESERVED to AVAILABLE
$ws.Range("E6").Value = "AVAILABLE"
$ws.Range("E7").Value = "AVAILABLE"
$ws.Range("E8").Value = "AVAILABLE"

# Update the current selection to match the edited range
$ws.Range("E5:E8").Select()
